$d = $word.ActiveDocument

# The "Contact Information" heading is the paragraph that should follow
# the new mission-statement paragraph. Get an insertion point collapsed
# to its very start.
$anchorPara = $d.Paragraphs(2)
$insertPoint = $d.Range($anchorPara.Range.Start, $anchorPara.Range.Start)

# Build the new paragraph as a WordprocessingML fragment (wrapped in the
# minimal flat-OPC "xmlPackage" envelope InsertXML expects) so it gets
# plain/body formatting -- no inherited Heading1 pStyle/direct
# formatting like InsertParagraphAfter would carry over, and no pStyle
# element at all (matching a normal/body paragraph). A trailing empty
# <w:p/> is required so the new paragraph actually gets its own
# paragraph mark instead of being merged into the following paragraph;
# that placeholder paragraph is removed immediately afterward.
$missionText = "Software engineer seeking opportunities to employ " + `
  "knowledge of web development to improve user experience in " + `
  "everyday tasks"

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:line="240" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
              </w:rPr>
              <w:t>$missionText</w:t>
            </w:r>
          </w:p>
          <w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xml) | Out-Null

# Drop the placeholder blank paragraph InsertXML had to materialize as
# the fragment's terminator.
$missionPara = $d.Paragraphs(2)
$placeholder = $missionPara.Next()
$placeholder.Range.Delete() | Out-Null
